$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "55.920.47"
$ws.Range("E2").Value = "  -0.91%  "
$ws.Range("D3").Value = "2.384.34"
$ws.Range("E3").Value = "  -4.90%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "478.68"
$ws.Range("E5").Value = "  -2.41%  "
$ws.Range("D6").Value = "147.34"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").Value = "0.501"
$ws.Range("E8").Value = "  -2.82%  "
$ws.Range("D9").Value = "2.388.78"
$ws.Range("E9").Value = "  -5.59%  "
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").Value = "5.46"
$ws.Range("E11").Value = "  -4.99%  "
$ws.Range("E12").Value = "  -3.32%  "
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "2.803.78"
$ws.Range("E14").Value = "  -4.62%  "
$ws.Range("D15").Value = "55.995.69"
$ws.Range("E15").Value = "  -0.69%  "
$ws.Range("D16").Value = "20.34"
$ws.Range("E16").Value = "  -4.36%  "
$ws.Range("E17").Value = "  -3.90%  "
$ws.Range("D18").Value = "2.379.81"
$ws.Range("E18").Value = "  -5.62%  "
$ws.Range("D19").Value = "4.53"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("D20").Value = "315.03"
$ws.Range("E20").Value = "  -2.48%  "
$ws.Range("D21").Value = "9.72"
$ws.Range("E21").Value = "  -5.46%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "5.69"
$ws.Range("E23").Value = "  -2.77%  "
$ws.Range("D24").Value = "56.74"
$ws.Range("E24").Value = "  -3.33%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.50%  "
$ws.Range("D26").Value = "0.395"
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("D27").Value = "0.157"
$ws.Range("E27").Value = "  -6.08%  "
$ws.Range("D28").Value = "2.497.70"
$ws.Range("E28").Value = "  -4.51%  "
$ws.Range("E29").Value = "  -5.17%  "
$ws.Range("D30").Value = "0.0₃0772"
$ws.Range("E30").Value = "  -3.64%  "
$ws.Range("E31").Value = "  +0.09%  "
$ws.Range("D32").Value = "147.52"
$ws.Range("E32").Value = "  -1.54%  "
$ws.Range("D33").Value = "18.02"
$ws.Range("E33").Value = "  -2.05%  "
$ws.Range("E34").Value = "  -1.94%  "
$ws.Range("D35").Value = "5.05"
$ws.Range("E35").Value = "  -3.49%  "
$ws.Range("E36").Value = "  -4.56%  "
$ws.Range("E37").Value = "  -4.65%  "
$ws.Range("D38").Value = "0.836"
$ws.Range("E38").Value = "  -4.59%  "
$ws.Range("D39").Value = "33.37"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("D42").Value = "3.39"
$ws.Range("E42").Value = "  -4.81%  "
$ws.Range("D43").Value = "0.0538"
$ws.Range("E43").Value = "  -3.72%  "
$ws.Range("D44").Value = "0.0946"
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("D45").Value = "0.584"
$ws.Range("E45").Value = "  -6.30%  "
$ws.Range("E46").Value = "  +0.41%  "
$ws.Range("D47").Value = "254.77"
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("D48").Value = "4.60"
$ws.Range("E48").Value = "  -5.58%  "
$ws.Range("E49").Value = "  -3.49%  "
$ws.Range("D50").Value = "17.01"
$ws.Range("E50").Value = "  -4.09%  "
$ws.Range("D51").Value = "1.780.28"
$ws.Range("E51").Value = "  -7.44%  "
